$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the formatting of the existing header row (B1)
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column H2:H15, all zeros (era data placeholder)
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
